# Repeat microarray experiments with min-max norm
#
# This script applies two related edits:
# 1. On the "experiment_plan" sheet, the cell-line identity values in rows
#    2-5 (the doppel-pair / repeat-experiment rows) are rotated between
#    two swapped pairs:
#       row 2: T84_LARGE_INTESTINE      -> HCT116_LARGE_INTESTINE
#       row 3: SNU283_LARGE_INTESTINE   -> SNUC4_LARGE_INTESTINE
#       row 4: SNUC4_LARGE_INTESTINE    -> SNU283_LARGE_INTESTINE
#       row 5: HCT116_LARGE_INTESTINE   -> T84_LARGE_INTESTINE
#    Only the "value" cells (the ones that currently hold the cell-line
#    name, not the class/batch id) are touched.
# 2. On the "train_valid_stats" sheet, the num_doppel_pairs count for the
#    Doppel_2 repeat drops from 3 to 2.

$wb = $excel.ActiveWorkbook

# ---- 1. experiment_plan sheet -------------------------------------------------
$plan = $wb.Worksheets.Item("experiment_plan")

# Each row has its own set of "value" (cell-line-name) columns - the
# remaining columns on that row hold the class/batch id and must not be
# touched. These were taken from the actual per-row cell-line columns.
$rowEdits = @(
    @{ Row = 2; Old = "T84_LARGE_INTESTINE";    New = "HCT116_LARGE_INTESTINE"; Cols = @("A","D","F","H","J","L") }
    @{ Row = 3; Old = "SNU283_LARGE_INTESTINE"; New = "SNUC4_LARGE_INTESTINE";  Cols = @("A","C","F","H","J","L") }
    @{ Row = 4; Old = "SNUC4_LARGE_INTESTINE";  New = "SNU283_LARGE_INTESTINE"; Cols = @("A","C","E","H","J","L") }
    @{ Row = 5; Old = "HCT116_LARGE_INTESTINE"; New = "T84_LARGE_INTESTINE";    Cols = @("A","C","E","G","J","L") }
)

foreach ($edit in $rowEdits) {
    $row = $edit.Row
    $old = $edit.Old
    $new = $edit.New
    foreach ($col in $edit.Cols) {
        $cell = $plan.Range("$col$row")
        if ($cell.Value2 -eq $old) {
            $cell.Value2 = $new
        }
    }
}

# ---- 2. train_valid_stats sheet -----------------------------------------------
$stats = $wb.Worksheets.Item("train_valid_stats")

# Find the row whose "Doppel" column (A) is Doppel_2, and set its
# num_doppel_pairs column (F) to 2.
for ($r = 2; $r -le 8; $r++) {
    $label = $stats.Cells.Item($r, 1).Value2
    if ($label -eq "Doppel_2") {
        $stats.Cells.Item($r, 6).Value2 = 2
        break
    }
}
